$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H8").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = 31
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = 5
$ws.Range("J10").Value = 16
$ws.Range("K10").Value = 102
$ws.Range("L10").Value = 73
$ws.Range("M10").Value = 9
$ws.Range("N10").Value = 68
$ws.Range("I13").Value = 5
$ws.Range("J13").Value = 661
$ws.Range("K13").Value = 10
$ws.Range("L13").Value = 910
$ws.Range("M13").Value = -544
$ws.Range("N13").Value = 544
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 3
$ws.Range("G14").Value = 7
$ws.Range("H14").Value = 50
$ws.Range("I14").Value = 73
$ws.Range("J14").Value = 59
$ws.Range("K14").Value = 100
$ws.Range("L14").Value = -99
$ws.Range("M14").Value = 0
$ws.Range("E15").Value = 21
$ws.Range("F15").Value = 22
$ws.Range("G15").Value = 277
$ws.Range("H15").Value = -69
$ws.Range("I15").Value = 154
$ws.Range("J15").Value = 201
$ws.Range("K15").Value = 330
$ws.Range("L15").Value = 476
$ws.Range("M15").Value = 473
$ws.Range("N15").Value = 619
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 619
$ws.Range("G16").Value = 463
$ws.Range("H16").Value = 486
$ws.Range("I16").Value = 531
$ws.Range("J16").Value = 1697
$ws.Range("K16").Value = 523
$ws.Range("L16").Value = 862
$ws.Range("M16").Value = 662
$ws.Range("N16").Value = 947
$ws.Range("E17").Value = 11102
$ws.Range("F17").Value = 15198
$ws.Range("G17").Value = 20858
$ws.Range("H17").Value = 17110
$ws.Range("I17").Value = 17386
$ws.Range("J17").Value = 21254
$ws.Range("K17").Value = 33321
$ws.Range("L17").Value = 25256
$ws.Range("M17").Value = 27610
$ws.Range("N17").Value = 27184
$ws.Range("E19").Value = 7871
$ws.Range("F19").Value = 37956
$ws.Range("G19").Value = 18545
$ws.Range("H19").Value = 24899
$ws.Range("I19").Value = 17004
$ws.Range("J19").Value = 36233
$ws.Range("K19").Value = 70140
$ws.Range("L19").Value = -22079
$ws.Range("M19").Value = 33173
$ws.Range("N19").Value = 13199
$ws.Range("E20").Value = 19308
$ws.Range("F20").Value = 53829
$ws.Range("G20").Value = 40150
$ws.Range("H20").Value = 42478
$ws.Range("I20").Value = 35158
$ws.Range("J20").Value = 60121
$ws.Range("K20").Value = 104526
$ws.Range("L20").Value = 5399
$ws.Range("M20").Value = 61383
$ws.Range("N20").Value = 42561
$ws.Range("E24").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F24").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G24").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H24").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I24").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J24").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K24").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L24").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M24").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N24").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("E26").Value = 149
$ws.Range("F26").Value = 160
$ws.Range("G26").Value = 149
$ws.Range("H26").Value = 155
$ws.Range("I26").Value = 150
$ws.Range("J26").Value = 151
$ws.Range("K26").Value = 145
$ws.Range("L26").Value = 147
$ws.Range("M26").Value = 143
$ws.Range("N26").Value = 148
$ws.Range("E27").Value = 92
$ws.Range("F27").Value = 83
$ws.Range("G27").Value = 96
$ws.Range("H27").Value = 102
$ws.Range("K27").Value = 102
$ws.Range("M27").Value = 101
$ws.Range("N27").Value = 102
